$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2427
$ws.Range("C2").Value = 57.39

$ws.Range("B3").Value = 1028
$ws.Range("C3").Value = 24.31

$ws.Range("B4").Value = 546
$ws.Range("C4").Value = 12.91

$ws.Range("B5").Value = 216
$ws.Range("C5").Value = 5.11

$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 0.28
